$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws3 = $wb.Worksheets.Item("Hoja3")

# Clear the header row (row 1) contents on Hoja3
$ws3.Range("A1:D1").ClearContents()

# Copy the "Twitter" data row from Hoja1 (row 3) into Hoja3 row 3, preserving
# the original cell types (text "FALSE" rather than a Boolean).
$ws1.Range("A3:D3").Copy()
$ws3.Range("A3").PasteSpecial(-4163)

# Select cell G5 and make Hoja3 the active sheet
$ws3.Activate()
$ws3.Range("G5").Select()
